$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to Text format so values that look numeric (e.g. "1.05")
    # are kept verbatim as strings instead of being parsed into numbers,
    # then restore the default "Normal" style so no stray formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "62.479.97"
$ws.Range("E2").Value = "  -2.41%  "
$ws.Range("D3").Value = "3.203.31"
$ws.Range("E3").Value = "  -3.42%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue $ws.Range("D5") "595.15"
$ws.Range("E5").Value = "  -0.76%  "
Set-TextValue $ws.Range("D6") "136.73"
$ws.Range("E6").Value = "  -4.63%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "3.202.14"
$ws.Range("E8").Value = "  -3.41%  "
$ws.Range("E9").Value = "  -2.89%  "
Set-TextValue $ws.Range("D10") "0.145"
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("E11").Value = "  -2.27%  "
Set-TextValue $ws.Range("D12") "0.456"
$ws.Range("E12").Value = "  -3.72%  "
$ws.Range("E13").Value = "  -4.84%  "
Set-TextValue $ws.Range("D14") "33.67"
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").Value = "3.727.07"
$ws.Range("E15").Value = "  -3.49%  "
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "3.199.59"
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("D18").Value = "62.594.80"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("E19").Value = "  -3.03%  "
Set-TextValue $ws.Range("D20") "463.10"
$ws.Range("E20").Value = "  -4.14%  "
$ws.Range("E21").Value = "  -2.23%  "
Set-TextValue $ws.Range("D22") "0.714"
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("E23").Value = "  -3.97%  "
Set-TextValue $ws.Range("D24") "13.57"
$ws.Range("E24").Value = "  -0.15%  "
Set-TextValue $ws.Range("D25") "83.94"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -2.47%  "
Set-TextValue $ws.Range("D28") "0.999"
$ws.Range("E28").Value = "  -0.15%  "
Set-TextValue $ws.Range("D29") "7.96"
$ws.Range("E29").Value = "  -3.37%  "
Set-TextValue $ws.Range("D30") "6.97"
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("E31").Value = "  -3.47%  "
Set-TextValue $ws.Range("D32") "27.46"
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("E33").Value = "  -3.49%  "
$ws.Range("E34").Value = "  -4.50%  "
Set-TextValue $ws.Range("D35") "1.05"
$ws.Range("E35").Value = "  -5.24%  "
Set-TextValue $ws.Range("D36") "5.89"
$ws.Range("E36").Value = "  -1.92%  "
Set-TextValue $ws.Range("D37") "51.65"
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("D38").Value = "0.0₃0697"
$ws.Range("E38").Value = "  -8.33%  "
Set-TextValue $ws.Range("D39") "0.0390"
$ws.Range("E39").Value = "  -2.98%  "
Set-TextValue $ws.Range("D40") "416.08"
$ws.Range("E40").Value = "  -4.12%  "
$ws.Range("D41").Value = "2.996.75"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("E42").Value = "  +3.82%  "
Set-TextValue $ws.Range("D43") "8.12"
$ws.Range("E43").Value = "  -3.79%  "
Set-TextValue $ws.Range("D44") "2.67"
$ws.Range("E44").Value = "  -4.82%  "
Set-TextValue $ws.Range("D45") "2.20"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("E46").Value = "  -6.14%  "
Set-TextValue $ws.Range("D47") "36.34"
$ws.Range("E47").Value = "  +2.87%  "
Set-TextValue $ws.Range("D48") "0.999"
$ws.Range("E48").Value = "  -0.05%  "
Set-TextValue $ws.Range("D49") "25.97"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D50") "124.13"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D51") "2.29"
$ws.Range("E51").Value = "  -1.76%  "
